$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 7: van den Berg et al. (2020) - AULs primer paper
# ------------------------------------------------------------------
$ws.Range("C7").Value = "https://link.springer.com/article/10.1007/s11214-020-00771-x"
$ws.Range("B7").Value = "A Primer on Focused Solar Energetic Particle Transport`nBasic Physics and Recent Modelling Results"
$ws.Range("A7").Value = "van den Berg"
$ws.Range("D7").Value = "https://doi.org/10.1007/s11214-020-00771-x"
$ws.Range("E7").Value = 2020

# ------------------------------------------------------------------
# Row 8: Steyn et al. (2020) - soft X-ray Neupert effect paper
# ------------------------------------------------------------------
$ws.Range("C8").Value = "https://www.swsc-journal.org/articles/swsc/full_html/2020/01/swsc200079/swsc200079.html"
$ws.Range("B8").Value = "The soft X-ray Neupert effect as a proxy for solar energetic particle injection`nA proof-of-concept physics-based forecasting model"
$ws.Range("D8").Value = "https://doi.org/10.1051/swsc/2020067"
$ws.Range("A8").Value = "Steyn"
$ws.Range("E8").Value = 2020

# ------------------------------------------------------------------
# Hyperlinks (matching column order hyperlinks were originally added)
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D7"), "https://doi.org/10.1007/s11214-020-00771-x")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://link.springer.com/article/10.1007/s11214-020-00771-x")
$ws.Hyperlinks.Add($ws.Range("C8"), "https://www.swsc-journal.org/articles/swsc/full_html/2020/01/swsc200079/swsc200079.html")
$ws.Hyperlinks.Add($ws.Range("D8"), "https://doi.org/10.1051/swsc/2020067")

# ------------------------------------------------------------------
# Formatting: copy the row-2/row-6 cell formats (wrap text, hyperlink
# style on link/doi columns) down onto the two new rows.
# ------------------------------------------------------------------
$ws.Range("A2:E2").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)

$ws.Range("A2:E2").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)

$ws.Rows.Item(7).RowHeight = 51
$ws.Rows.Item(8).RowHeight = 68

# ------------------------------------------------------------------
# Selection / view state tweak (matches the authored workbook view)
# ------------------------------------------------------------------
$ws.Range("D15").Select()
